# Generate Report for Handback
# Refresh the generated timestamps recorded on the handback-status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet - Latest HO Xliff Generate Date for first file row.
# This timestamp is shared (same text) with the de-de sheet's
# "Correspond Handoff Datetime" for the same file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 13:07:59"

# "zh-cn" sheet - Correspond Handoff / Handback Datetime for first file row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-21 13:07:55"
$wsZhCn.Range("K2").Value = "2016-08-21 13:08:14"

# "de-de" sheet - Correspond Handoff / Handback Datetime for first file row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-21 13:07:59"
$wsDeDe.Range("K2").Value = "2016-08-21 13:08:20"
